$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date refresh
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# The old duplicate "Contact / No display for ContactDetail" rows (10 & 11)
# become a single "Jurisdiction / United States of America" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-redundant duplicate row (old row 11), shifting everything up.
$ws.Rows.Item(11).Delete()

# --- Elements sheet updates ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition now describe this specific extension
# instead of the generic "Extension" / "An Extension" placeholders.
$ws2.Range("K2").Value = "Vision Fully Insured Indicator"
$ws2.Range("L2").Value = "Indicator of the fully insured vision coverage for the member or employee"
